$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.384.81"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.800.61"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.24"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.44%  "
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0677"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0964"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").Value = "2.059.51"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.28"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").Value = "1.808.41"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.627"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D16").Value = "34.384.11"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.42"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.12"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("E20").Value = "  -2.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.22"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("E23").Value = "  -1.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.76"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  +4.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.33"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.69%  "
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("E30").Value = "  -1.47%  "
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("E32").Value = "  -1.40%  "
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.78"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.03%  "
$ws.Range("D35").Value = "1.361.45"
$ws.Range("E35").Value = "  -2.55%  "
$ws.Range("E36").Value = "  -3.86%  "
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.35"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -6.37%  "
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("E40").Value = "  +1.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.77"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "80.51"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.936"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("E44").Value = "  +5.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.26"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("D47").Value = "1.962.80"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  -3.44%  "
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.25"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("E51").Value = "  -3.01%  "
